$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.374.27"
$ws.Range("E2").Value = "  +1.95%  "
$ws.Range("D3").Value = "2.596.17"
$ws.Range("E3").Value = "  +0.51%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "568.76"
$ws.Range("E5").Value = "  +1.47%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.97"
$ws.Range("E6").Value = "  -0.67%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.997"
$ws.Range("E7").Value = "  -0.17%  "
$ws.Range("E8").Value = "  +0.53%  "
$ws.Range("D9").Value = "2.616.94"
$ws.Range("E9").Value = "  +1.02%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.55"
$ws.Range("E10").Value = "  -1.49%  "
$ws.Range("E11").Value = "  +1.18%  "
$ws.Range("E12").Value = "  +2.85%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.150"
$ws.Range("E13").Value = "  -5.86%  "
$ws.Range("D14").Value = "3.063.21"
$ws.Range("E14").Value = "  +0.78%  "
$ws.Range("D15").Value = "60.378.32"
$ws.Range("E15").Value = "  +2.01%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "23.27"
$ws.Range("E16").Value = "  +1.19%  "
$ws.Range("E17").Value = "  +2.84%  "
$ws.Range("D18").Value = "2.618.66"
$ws.Range("E18").Value = "  +1.53%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.32"
$ws.Range("E19").Value = "  +9.38%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.65"
$ws.Range("E20").Value = "  +1.63%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "346.03"
$ws.Range("E21").Value = "  +2.79%  "
$ws.Range("E22").Value = "  +8.60%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.998"
$ws.Range("E23").Value = "  -0.29%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.528"
$ws.Range("E24").Value = "  +13.33%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "63.20"
$ws.Range("E25").Value = "  -1.27%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  +0.06%  "
$ws.Range("E27").Value = "  -1.75%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.68"
$ws.Range("E28").Value = "  +4.60%  "
$ws.Range("D29").Value = "0.0₃0783"
$ws.Range("E29").Value = "  +1.11%  "
$ws.Range("E30").Value = "  +9.08%  "
$ws.Range("E31").Value = "  +4.03%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.997"
$ws.Range("E32").Value = "  -0.08%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "160.35"
$ws.Range("E33").Value = "  +0.94%  "
$ws.Range("E34").Value = "  +2.49%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.22"
$ws.Range("E35").Value = "  +4.78%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.966"
$ws.Range("E36").Value = "  +9.81%  "
$ws.Range("E37").Value = "  +3.85%  "
$ws.Range("E38").Value = "  +8.63%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "37.72"
$ws.Range("E39").Value = "  +0.67%  "
$ws.Range("E40").Value = "  +3.59%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.847"
$ws.Range("E41").Value = "  -2.34%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "294.23"
$ws.Range("E42").Value = "  +0.60%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "138.80"
$ws.Range("E43").Value = "  +4.45%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.996"
$ws.Range("E44").Value = "  -0.33%  "
$ws.Range("E45").Value = "  +0.87%  "
$ws.Range("E46").Value = "  +1.58%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "19.68"
$ws.Range("E47").Value = "  +3.51%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0544"
$ws.Range("E48").Value = "  +1.86%  "
$ws.Range("E49").Value = "  +2.68%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "19.81"
$ws.Range("E50").Value = "  +6.54%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "10.71"
$ws.Range("E51").Value = "  +0.76%  "
